$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for each data row (rows 2-170).
# All of these cells currently equal 2023-09-06 (serial 45175) and must be updated
# to 2023-09-08 (serial 45177), while keeping their existing date formatting/style.
$oldSerial = 45175
$newSerial = 45177

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 170 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldSerial) {
        $cell.Value = $newSerial
    }
}
